# repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows whose final-step deltas
# were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 1
    7  = -2
    8  = -5
    9  = -1
    10 = 1
    13 = -3
    14 = -2
    18 = 0
    21 = 1
    23 = 7
    31 = 0
    35 = 9
    40 = -2
    43 = 1
    45 = 0
    50 = 2
    53 = 1
    57 = 3
    63 = 3
    64 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
